$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3: Gabriel/50 -> Clebinho/40
$ws.Range("A3").Value = "Clebinho"
$ws.Range("B3").Value = 40

# Update row 4: Clebinho/50 -> Teste/0
$ws.Range("A4").Value = "Teste"
$ws.Range("B4").Value = 0

# Remove row 5 (Diego/50) entirely - delete the whole row so data ends at row 4
$ws.Rows("5").Delete()
